$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole edited block to Text format first so that numeric-
# looking strings (prices, percentages) are stored as text, matching
# the original inline-string cell content instead of being auto-converted
# to numbers by Excel's type inference. ClearFormats afterwards restores
# the original (default) cell style so no stray formatting is left behind.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "293.27"
$ws.Range("E2").Value = "0.10%"
$ws.Range("D3").Value = "40.39"
$ws.Range("E3").Value = "0.47%"
$ws.Range("D4").Value = "5.004"
$ws.Range("E4").Value = "-0.40%"
$ws.Range("D5").Value = "0.07342"
$ws.Range("E5").Value = "-0.65%"
$ws.Range("D6").Value = "4.285"
$ws.Range("E6").Value = "-0.79%"
$ws.Range("D7").Value = "1.561"
$ws.Range("E7").Value = "2.81%"
$ws.Range("D8").Value = "0.9240"
$ws.Range("E8").Value = "0.14%"
$ws.Range("D9").Value = "2.345"
$ws.Range("E9").Value = "-2.25%"
$ws.Range("D10").Value = "0.1178"
$ws.Range("E10").Value = "1.07%"
$ws.Range("D11").Value = "0.1815"
$ws.Range("E11").Value = "3.53%"
$ws.Range("D12").Value = "0.04383"
$ws.Range("E12").Value = "4.85%"
$ws.Range("D13").Value = "0.08792"
$ws.Range("E13").Value = "1.77%"
$ws.Range("D14").Value = "0.1054"
$ws.Range("E14").Value = "0.09%"
$ws.Range("D15").Value = "0.001269"
$ws.Range("E15").Value = "-0.46%"
$ws.Range("D16").Value = "0.006001"
$ws.Range("E16").Value = "2.00%"
$ws.Range("E17").Value = "-0.59%"
$ws.Range("D18").Value = "0.3305"
$ws.Range("E18").Value = "-0.28%"
$ws.Range("D19").Value = "7.829"
$ws.Range("E19").Value = "3.31%"
$ws.Range("D20").Value = "0.1390"
$ws.Range("E20").Value = "2.35%"
$ws.Range("D21").Value = "0.2799"
$ws.Range("E21").Value = "-0.63%"
$ws.Range("D22").Value = "0.03922"
$ws.Range("E22").Value = "2.31%"
$ws.Range("D23").Value = "0.001261"
$ws.Range("E23").Value = "-1.87%"
$ws.Range("E24").Value = "1.50%"
$ws.Range("E25").Value = "-8.14%"
$ws.Range("D26").Value = "0.0003722"
$ws.Range("E26").Value = "-0.48%"
$ws.Range("D38").Value = "0.02341"
$ws.Range("E38").Value = "1.08%"
$ws.Range("D39").Value = "0.05105"
$ws.Range("E39").Value = "2.18%"
$ws.Range("D41").Value = "0.007865"
$ws.Range("E41").Value = "1.96%"
$ws.Range("D42").Value = "0.1292"
$ws.Range("E42").Value = "1.36%"
$ws.Range("D43").Value = "0.007377"
$ws.Range("E43").Value = "-0.63%"
$ws.Range("D44").Value = "0.008038"
$ws.Range("E44").Value = "1.62%"
$ws.Range("D45").Value = "0.2919"
$ws.Range("E45").Value = "-7.79%"
$ws.Range("D46").Value = "0.00006229"
$ws.Range("E46").Value = "-3.97%"
$ws.Range("E47").Value = "-0.47%"
$ws.Range("D48").Value = "0.04843"
$ws.Range("E48").Value = "-80.77%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.004200"
$ws.Range("E49").Value = "-0.48%"
$ws.Range("B50").Value = "CryptobidCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.47%"
$ws.Range("B51").Value = "SpecialPowerGold"
$ws.Range("C51").Value = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.47%"

$editRange.ClearFormats()

